# Apply updated cryptocurrency price/volume data to the worksheet.
# Source: automated "cryptos list" refresh (GitHub Actions).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new value is a plain decimal number (e.g. "533.26") must be
# forced to Text format first, otherwise Excel auto-converts the assigned
# string into a floating point number (losing the original text formatting
# and introducing binary floating point rounding). The number format is
# reset back to the default ("Normal" style) immediately afterwards so the
# cell keeps no extra formatting, matching the original inline-string cells.
function Set-TextValue($cell, $value) {
    $rng = $ws.Range($cell)
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.Style = "Normal"
}

# Row 2
$ws.Range('D2').Value = '58.493.81'
$ws.Range('E2').Value = '  -1.64%  '

# Row 3
$ws.Range('D3').Value = '2.620.70'
$ws.Range('E3').Value = '  +0.73%  '

# Row 4
$ws.Range('E4').Value = '  +0.02%  '

# Row 5
Set-TextValue 'D5' '533.26'
$ws.Range('E5').Value = '  -0.85%  '

# Row 6
Set-TextValue 'D6' '142.41'

# Row 7
$ws.Range('E7').Value = '  -0.01%  '

# Row 8
Set-TextValue 'D8' '0.568'
$ws.Range('E8').Value = '  +0.41%  '

# Row 9
$ws.Range('E9').Value = '  +6.48%  '

# Row 10
$ws.Range('E10').Value = '  -2.03%  '

# Row 11
$ws.Range('E11').Value = '  -0.44%  '

# Row 12
$ws.Range('E12').Value = '  +1.01%  '

# Row 13
$ws.Range('D13').Value = '3.085.97'
$ws.Range('E13').Value = '  +0.92%  '

# Row 14
$ws.Range('D14').Value = '58.418.08'
$ws.Range('E14').Value = '  -1.56%  '

# Row 15
$ws.Range('E15').Value = '  -0.44%  '

# Row 16
$ws.Range('D16').Value = '2.635.84'
$ws.Range('E16').Value = '  +1.97%  '

# Row 17
$ws.Range('E17').Value = '  -1.36%  '

# Row 18
Set-TextValue 'D18' '4.39'
$ws.Range('E18').Value = '  +0.35%  '

# Row 19
Set-TextValue 'D19' '334.19'
$ws.Range('E19').Value = '  -2.08%  '

# Row 20
Set-TextValue 'D20' '10.11'
$ws.Range('E20').Value = '  +0.27%  '

# Row 21
Set-TextValue 'D21' '6.22'
$ws.Range('E21').Value = '  -2.28%  '

# Row 22
$ws.Range('E22').Value = '  +0.02%  '

# Row 23
Set-TextValue 'D23' '66.38'
$ws.Range('E23').Value = '  -1.47%  '

# Row 24
$ws.Range('E24').Value = '  +1.43%  '

# Row 25
$ws.Range('B25').Value = 'Kaspa'
$ws.Range('C25').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
Set-TextValue 'D25' '0.163'
$ws.Range('E25').Value = '  -1.14%  '

# Row 26
$ws.Range('B26').Value = 'Binance-PegBSC-USD'
$ws.Range('C26').Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
Set-TextValue 'D26' '0.999'
$ws.Range('E26').Value = '  +0.04%  '

# Row 27
Set-TextValue 'D27' '7.07'
$ws.Range('E27').Value = '  -2.12%  '

# Row 28
Set-TextValue 'D28' '0.998'
$ws.Range('E28').Value = '  -0.09%  '

# Row 29
$ws.Range('D29').Value = '0.0₃0732'
$ws.Range('E29').Value = '  -1.58%  '

# Row 30
$ws.Range('E30').Value = '  -1.47%  '

# Row 31
$ws.Range('E31').Value = '  +0.32%  '

# Row 32
Set-TextValue 'D32' '18.74'
$ws.Range('E32').Value = '  -0.41%  '

# Row 33
Set-TextValue 'D33' '150.17'
$ws.Range('E33').Value = '  +0.24%  '

# Row 34
$ws.Range('E34').Value = '  -2.44%  '

# Row 35
Set-TextValue 'D35' '0.853'
$ws.Range('E35').Value = '  +1.31%  '

# Row 36
Set-TextValue 'D36' '1.10'
$ws.Range('E36').Value = '  -1.62%  '

# Row 37
$ws.Range('E37').Value = '  -3.58%  '

# Row 38
Set-TextValue 'D38' '0.808'

# Row 39
Set-TextValue 'D39' '3.56'
$ws.Range('E39').Value = '  +0.76%  '

# Row 40
Set-TextValue 'D40' '279.67'
$ws.Range('E40').Value = '  +2.94%  '

# Row 41
$ws.Range('E41').Value = '  +0.03%  '

# Row 42
Set-TextValue 'D42' '0.593'
$ws.Range('E42').Value = '  -1.10%  '

# Row 43
Set-TextValue 'D43' '10.66'
$ws.Range('E43').Value = '  -0.61%  '

# Row 44
$ws.Range('B44').Value = 'EnergySwap'
$ws.Range('C44').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
Set-TextValue 'D44' '18.96'
$ws.Range('E44').Value = '  +2.41%  '

# Row 45
$ws.Range('B45').Value = 'Hedera'
$ws.Range('C45').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
Set-TextValue 'D45' '0.0527'
$ws.Range('E45').Value = '  +0.58%  '

# Row 46
Set-TextValue 'D46' '0.0935'
$ws.Range('E46').Value = '  -1.87%  '

# Row 47
$ws.Range('E47').Value = '  +0.19%  '

# Row 48
$ws.Range('D48').Value = '1.937.14'
$ws.Range('E48').Value = '  -0.23%  '

# Row 49
Set-TextValue 'D49' '4.42'
$ws.Range('E49').Value = '  -1.74%  '

# Row 50
Set-TextValue 'D50' '17.83'
$ws.Range('E50').Value = '  -4.14%  '

# Row 51
Set-TextValue 'D51' '112.87'
$ws.Range('E51').Value = '  +0.80%  '
